$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data captured from the target commit, keyed by destination row.
# Each row pulls its new Fecha/Volumen/Precio* values from another row
# in the original (pre-edit) snapshot -- see analysis in the task.
$rows = @(
    [pscustomobject]@{ Row = 2; D = 44957; J = 20; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 3; D = 44959; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 4; D = 44781; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 6; D = 44498; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 7; D = 45194; J = 40; K = 6000; L = 6000; M = 6000; P = 6000 }
    [pscustomobject]@{ Row = 8; D = 44749; J = 65; K = 6000; L = 6000; M = 6000; P = 6000 }
    [pscustomobject]@{ Row = 9; D = 45169; J = 50; K = 4000; L = 5000; M = 4600; P = 4600 }
    [pscustomobject]@{ Row = 10; D = 44312; J = 50; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 11; D = 44291; J = 35; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 13; D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
    [pscustomobject]@{ Row = 14; D = 44966; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 15; D = 44504; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 16; D = 44259; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 17; D = 44315; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 18; D = 44497; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 20; D = 45163; J = 30; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 21; D = 44656; J = 85; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 23; D = 44679; J = 50; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 24; D = 44176; J = 10; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 25; D = 44316; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 26; D = 45159; J = 75; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 27; D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 28; D = 44509; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 29; D = 44680; J = 20; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 30; D = 44956; J = 40; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 31; D = 44508; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
    [pscustomobject]@{ Row = 32; D = 44365; J = 55; K = 5000; L = 5000; M = 5000; P = 5000 }
    [pscustomobject]@{ Row = 33; D = 44649; J = 20; K = 5000; L = 5000; M = 5000; P = 5000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
